$d = $word.ActiveDocument

# 1. First paragraph in the body: spacing before/after 240 (12pt) -> 300 (15pt)
$p1 = $d.Paragraphs(1)
$p1.Format.SpaceBefore = 15
$p1.Format.SpaceAfter = 15

# 2. Header: paragraph alignment right -> left
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
for ($i = 1; $i -le $hdr.Range.Paragraphs.Count; $i++) {
    $hdr.Range.Paragraphs($i).Format.Alignment = 0
}

# 3. Footer: remove the first paragraph (PAGE field, right aligned) entirely,
#    and set the remaining (previously centered) paragraph's alignment to left.
$ftr = $sec.Footers(1)
$ftr.Range.Paragraphs(1).Range.Delete()
for ($i = 1; $i -le $ftr.Range.Paragraphs.Count; $i++) {
    $ftr.Range.Paragraphs($i).Format.Alignment = 0
}
